# Insert a new row for the missing "get_reward" support_module_effect entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 70, pushing the existing rows (old 70+) down by one.
$ws.Rows.Item(70).Insert()

# Fill the new row 70 with the support_module_effect / get_reward entry.
# (Column E is written before column B so the new shared-string entries
# land in the same order as the target workbook: "보상 획득" then "get_reward".)
$ws.Cells.Item(70, 1).Value = "support_module_effect"
$ws.Cells.Item(70, 5).Value = "보상 획득"
$ws.Cells.Item(70, 2).Value = "get_reward"
$ws.Cells.Item(70, 3).Formula = '=UPPER(A70)&"_"&UPPER(B70)'
$ws.Cells.Item(70, 4).Value = 5

# Restore the view state to match the post-edit workbook.
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("C70").Select()
